{"js": "const replacements = [\n  [\"160\u00f74=40, 0\", \"240\u00f74=60, 0\"],\n  [\"642\u00f74=160, 2\", \"317\u00f72=158, 1\"],\n  [\"657\u00f77=93, 6\", \"224\u00f75=44, 4\"],\n  [\"569\u00f77=81, 2\", \"450\u00f72=225, 0\"],\n  [\"155\u00f76=25, 5\", \"821\u00f76=136, 5\"],\n  [\"778\u00f79=86, 4\", \"112\u00f72=56, 0\"],\n  [\"654\u00f76=109, 0\", \"711\u00f72=355, 1\"],\n  [\"858\u00f75=171, 3\", \"254\u00f75=50, 4\"],\n  [\"904\u00f72=452, 0\", \"469\u00f76=78, 1\"],\n  [\"743\u00f76=123, 5\", \"570\u00f77=81, 3\"],\n  [\"777\u00f79=86, 3\", \"341\u00f72=170, 1\"],\n  [\"580\u00f74=145, 0\", \"618\u00f75=123, 3\"],\n  [\"915\u00f72=457, 1\", \"341\u00f78=42, 5\"],\n  [\"662\u00f74=165, 2\", \"558\u00f75=111, 3\"],\n  [\"483\u00f73=161, 0\", \"997\u00f72=498, 1\"],\n  [\"443\u00f72=221, 1\", \"815\u00f72=407, 1\"],\n  [\"398\u00f76=66, 2\", \"959\u00f79=106, 5\"],\n  [\"261\u00f77=37, 2\", \"995\u00f76=165, 5\"],\n  [\"606\u00f74=151, 2\", \"153\u00f75=30, 3\"],\n  [\"285\u00f78=35, 5\", \"343\u00f73=114, 1\"],\n  [\"310\u00f76=51, 4\", \"589\u00f78=73, 5\"],\n  [\"491\u00f77=70, 1\", \"737\u00f78=92, 1\"],\n  [\"922\u00f76=153, 4\", \"856\u00f72=428, 0\"],\n  [\"529\u00f74=132, 1\", \"252\u00f77=36, 0\"],\n  [\"469\u00f72=234, 1\", \"352\u00f79=39, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Not found: ${oldText}`);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"160\u00f74=40, 0\", \"240\u00f74=60, 0\"),\n    @(\"642\u00f74=160, 2\", \"317\u00f72=158, 1\"),\n    @(\"657\u00f77=93, 6\", \"224\u00f75=44, 4\"),\n    @(\"569\u00f77=81, 2\", \"450\u00f72=225, 0\"),\n    @(\"155\u00f76=25, 5\", \"821\u00f76=136, 5\"),\n    @(\"778\u00f79=86, 4\", \"112\u00f72=56, 0\"),\n    @(\"654\u00f76=109, 0\", \"711\u00f72=355, 1\"),\n    @(\"858\u00f75=171, 3\", \"254\u00f75=50, 4\"),\n    @(\"904\u00f72=452, 0\", \"469\u00f76=78, 1\"),\n    @(\"743\u00f76=123, 5\", \"570\u00f77=81, 3\"),\n    @(\"777\u00f79=86, 3\", \"341\u00f72=170, 1\"),\n    @(\"580\u00f74=145, 0\", \"618\u00f75=123, 3\"),\n    @(\"915\u00f72=457, 1\", \"341\u00f78=42, 5\"),\n    @(\"662\u00f74=165, 2\", \"558\u00f75=111, 3\"),\n    @(\"483\u00f73=161, 0\", \"997\u00f72=498, 1\"),\n    @(\"443\u00f72=221, 1\", \"815\u00f72=407, 1\"),\n    @(\"398\u00f76=66, 2\", \"959\u00f79=106, 5\"),\n    @(\"261\u00f77=37, 2\", \"995\u00f76=165, 5\"),\n    @(\"606\u00f74=151, 2\", \"153\u00f75=30, 3\"),\n    @(\"285\u00f78=35, 5\", \"343\u00f73=114, 1\"),\n    @(\"310\u00f76=51, 4\", \"589\u00f78=73, 5\"),\n    @(\"491\u00f77=70, 1\", \"737\u00f78=92, 1\"),\n    @(\"922\u00f76=153, 4\", \"856\u00f72=428, 0\"),\n    @(\"529\u00f74=132, 1\", \"252\u00f77=36, 0\"),\n    @(\"469\u00f72=234, 1\", \"352\u00f79=39, 1\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
